$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of treasury delta data to append (fix for fedrollover bug)
$newRows = @(
    @(20082100, 1630312000000, 1621385000000, -8927000000),
    @(20082200, 0, 0, 0),
    @(20082300, 0, 0, 0),
    @(20082400, 1621385000000, 1627631000000, 6246000000)
)

$startRow = 388
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
